$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 7.96
$ws.Range("C3").Value = 0

# Prepare formatting for the new rows in column A by copying the
# existing "index" style (bold font, border, centered) from A2.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add new rows of data
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "water"
$ws.Range("C4").Value = 840

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "methane"
$ws.Range("C5").Value = 400
